$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 18 (Price/Stock change) ---
$ws.Cells.Item(18,8).Value = 2500
$ws.Cells.Item(18,9).Value = 50

# --- Append new row 45 ---
$ws.Cells.Item(45,1).Value = "'nihal"
$ws.Cells.Item(45,2).Value = "'kamat"
$ws.Cells.Item(45,3).Value = "'Nihal Kamat"
$ws.Cells.Item(45,4).Value = "'9405920393"
$ws.Cells.Item(45,5).Value = "'Ponda"
$ws.Cells.Item(45,6).Value = "'Goa"
$ws.Cells.Item(45,7).Value = "'Kamat Stores"
$ws.Cells.Item(45,8).Value = 3500
$ws.Cells.Item(45,9).Value = 50
$ws.Cells.Item(45,10).Value = "'Available"
$ws.Cells.Item(45,11).Value = "'15.4027241187136"
$ws.Cells.Item(45,12).Value = "'74.0032517910004"

# Give the numeric Price/Stock cells the same style (quote-prefix cell style s="1")
# used by the rest of the row, by copying formatting from a cell that already
# carries it.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("H45:I45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Select the newly added row, mirroring the end-user's on-screen selection
$ws.Rows(45).Select() | Out-Null

Write-Host "done"
